# Applies the "WIP heating system simplification and energy calculations" edit:
#  - District heating sheet: rename "Warmtenet" -> "LT-Warmtenet" (row 5), add a new
#    "MT-Warmtenet" row (row 6) that duplicates row 5's shape with slightly different
#    numbers, and add a new "Verlies" (loss) column N with per-row loss fractions.
#  - Inhouse heat systems sheet: insert a new "Electric boiler" row above the
#    HT-radiator row, with a CAPEX computed as the average of two quoted prices.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "District heating" sheet (2nd tab)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# 1a. Rename existing "Warmtenet" entry (row 5) to "LT-Warmtenet"
$ws2.Cells.Item(5, 1).Value2 = "LT-Warmtenet"

# 1b. New row 6 - "MT-Warmtenet", a near-duplicate of row 5 (same units/format),
#     with an updated connection count (271 instead of 270).
$ws2.Range("A5:M5").Copy($ws2.Range("A6:M6"))
$ws2.Cells.Item(6, 1).Value2 = "MT-Warmtenet"
$ws2.Cells.Item(6, 2).Value2 = 271
$ws2.Cells.Item(6, 5).Value2 = 271
$ws2.Cells.Item(6, 13).Formula = "=L6*3.6/1000"

# 1c. New "Verlies" (loss) column N, with a header + one value per data row
$ws2.Cells.Item(1, 14).Value2 = "Verlies"
$ws2.Cells.Item(1, 14).Font.Bold = $true

$ws2.Cells.Item(2, 14).Value2 = 0.3
$ws2.Cells.Item(3, 14).Value2 = 0.3
$ws2.Cells.Item(4, 14).Value2 = 0.3
$ws2.Cells.Item(5, 14).Value2 = 0.15
$ws2.Cells.Item(6, 14).Value2 = 0.3

# View: scroll/select near the new column
$ws2.Range("M18").Select()

# ---------------------------------------------------------------------------
# 2. "Inhouse heat systems" sheet (1st tab) - insert "Electric boiler" row
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows.Item(7).Insert()

$ws1.Cells.Item(7, 1).Value2 = "Electric boiler"
$ws1.Cells.Item(7, 2).Value2 = "Decentral"
$ws1.Cells.Item(7, 3).Value2 = 0
$ws1.Cells.Item(7, 4).Value2 = "EUR/kW/y"
$ws1.Cells.Item(7, 5).Value2 = 0
$ws1.Cells.Item(7, 6).Formula = "=(1774.85+1112)/2"
$ws1.Cells.Item(7, 7).Value2 = "EUR per aansluiting"
$ws1.Cells.Item(7, 8).Value2 = 15

# View: restore the active sheet + selection expected after the edit
$ws1.Range("F7").Select()
